# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md file is now "Ready for handoff"
# (rather than "Handed back: in sync with en-US"), records the new
# handoff timestamps, and records an Error Detail explaining that the
# handback file version is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/429c14039e2ed7bf7b7f95a19f7eb8a1462d079b/e2e/a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45a52a35489bb18572efee076039d115692b7ecc/e2e/a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md."

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-22 16:49:22"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-22 16:49:17"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.140625

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-22 16:49:22"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.140625
